# Update the lx experiment results with new computed metrics.
# Only the numeric results for the "random_forest" (row 2) and "lsboost"
# (row 3) models changed; the header row, row labels, and the "old_model"
# (row 4) row keep their original values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# random_forest row (row 2): RMSE, NRMSE, MAE, RSE, RRSE, RAE, R2, Corr Coeff
$ws.Range("B2").Value = 2.1434986266114011
$ws.Range("C2").Value = 0.1873687610674302
$ws.Range("D2").Value = 1.5105784391534391
$ws.Range("E2").Value = 0.34253704223276443
$ws.Range("F2").Value = 0.58526664199556466
$ws.Range("G2").Value = 0.53339634150898296
$ws.Range("H2").Value = 0.65746295776723551
$ws.Range("I2").Value = 0.81395372334597615

# lsboost row (row 3): RMSE, NRMSE, MAE, RSE, RRSE, RAE, R2, Corr Coeff
$ws.Range("B3").Value = 1.4648382807338998
$ws.Range("C3").Value = 0.12804530425995631
$ws.Range("D3").Value = 1.0402304647540634
$ws.Range("E3").Value = 0.15997071500912871
$ws.Range("F3").Value = 0.39996339208623671
$ws.Range("G3").Value = 0.36731301721541798
$ws.Range("H3").Value = 0.84002928499087126
$ws.Range("I3").Value = 0.92550664357202717
